$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.736.00'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").Value = '1.627.35'
$ws.Range("E3").Value = '  -0.16%  '

$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = '  -0.74%  '

$ws.Range("D5").Value = "'214.43"
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("D6").Value = "'0.501"
$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("E7").Value = '  -0.60%  '

$ws.Range("E8").Value = '  -0.82%  '

$ws.Range("E9").Value = '  -0.60%  '

$ws.Range("D10").Value = "'19.56"
$ws.Range("E10").Value = '  +0.58%  '

$ws.Range("D11").Value = "'0.0792"
$ws.Range("E11").Value = '  +1.05%  '

$ws.Range("D12").Value = "'4.25"
$ws.Range("E12").Value = '  +0.30%  '

$ws.Range("D13").Value = '1.852.90'
$ws.Range("E13").Value = '  -0.10%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.633.21'
$ws.Range("E14").Value = '  +0.29%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = "'0.552"
$ws.Range("E15").Value = '  +0.44%  '

$ws.Range("D16").Value = '0.0₃0761'
$ws.Range("E16").Value = '  -0.83%  '

$ws.Range("D17").Value = "'62.78"
$ws.Range("E17").Value = '  -0.55%  '

$ws.Range("D18").Value = '25.725.49'
$ws.Range("E18").Value = '  -0.04%  '

$ws.Range("E19").Value = '  -0.59%  '

$ws.Range("D20").Value = "'4.44"
$ws.Range("E20").Value = '  +0.34%  '

$ws.Range("D21").Value = "'190.87"
$ws.Range("E21").Value = '  -1.46%  '

$ws.Range("E22").Value = '  -0.31%  '

$ws.Range("D23").Value = "'6.26"
$ws.Range("E23").Value = '  +0.95%  '

$ws.Range("E24").Value = '  -0.70%  '

$ws.Range("D25").Value = "'1.81"
$ws.Range("E25").Value = '  +1.67%  '

$ws.Range("D26").Value = "'142.53"
$ws.Range("E26").Value = '  +1.48%  '

$ws.Range("E27").Value = '  +3.28%  '

$ws.Range("E28").Value = '  +0.47%  '

$ws.Range("D29").Value = "'15.47"
$ws.Range("E29").Value = '  -0.13%  '

$ws.Range("E30").Value = '  +0.12%  '

$ws.Range("E31").Value = '  +1.88%  '

$ws.Range("E32").Value = '  -0.46%  '

$ws.Range("E33").Value = '  -0.73%  '

$ws.Range("E34").Value = '  -0.14%  '

$ws.Range("E35").Value = '  -0.25%  '

$ws.Range("E36").Value = '  +1.39%  '

$ws.Range("D37").Value = '1.136.98'
$ws.Range("E37").Value = '  +2.96%  '

$ws.Range("D38").Value = "'2.51"
$ws.Range("E38").Value = '  -2.03%  '

$ws.Range("E39").Value = '  -0.93%  '

$ws.Range("E40").Value = '  +0.03%  '

$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = '  -0.64%  '

$ws.Range("E42").Value = '  -1.02%  '

$ws.Range("D43").Value = "'5.57"
$ws.Range("E43").Value = '  -0.10%  '

$ws.Range("D44").Value = "'100.69"
$ws.Range("E44").Value = '  +0.84%  '

$ws.Range("E45").Value = '  +0.41%  '

$ws.Range("D46").Value = '1.762.66'

$ws.Range("D47").Value = "'55.11"
$ws.Range("E47").Value = '  +0.26%  '

$ws.Range("E48").Value = '  +1.71%  '

$ws.Range("D49").Value = "'1.45"
$ws.Range("E49").Value = '  +5.91%  '

$ws.Range("D50").Value = "'0.416"
$ws.Range("E50").Value = '  -0.63%  '

$ws.Range("D51").Value = "'2.33"
$ws.Range("E51").Value = '  +0.30%  '

$ws.Range("D4,D5,D6,D10,D11,D12,D15,D17,D20,D21,D23,D25,D26,D29,D38,D41,D43,D44,D47,D49,D50,D51").ClearFormats()
